# Generate Report for Handoff
# Adds a new "Ready for handoff" row for file
# 19ad70f1-efcd-41f1-9f5d-af48289612aeoooo....md
# to the Overview, zh-cn and de-de sheets/tables.

$wb = $excel.ActiveWorkbook

$fileName    = "19ad70f1-efcd-41f1-9f5d-af48289612aeooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$pathAndName = "e2e\19ad70f1-efcd-41f1-9f5d-af48289612aeooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$extension   = ".md"
$status      = "Ready for handoff"
$hoDate      = "2016-09-02 22:32:02"
$hoDateZh    = "2016-09-02 22:31:57"
$sourcePath  = "e2e"
$priority    = "ht"
$contentDup  = "False"
$zhXlf       = "19ad70f1-efcd-41f1-9f5d-af48289612aeoooooooooooooooooooooooooooooooooooooooo.ed99a82be1c5c7e59d3f190c55df76c473f30d1d.zh-cn.xlf"
$deXlf       = "19ad70f1-efcd-41f1-9f5d-af48289612aeoooooooooooooooooooooooooooooooooooooooo.ed99a82be1c5c7e59d3f190c55df76c473f30d1d.de-de.xlf"
$handbackDt  = "0001-01-01 00:00:00"
$toLocalize  = "True"
$hasMeta     = "False"

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/79b976807a68c4df918fee4e117ba6b6bb43deca/e2e/" + $fileName

# ---------------------------------------------------------------
# Sheet "Overview": File Name | Path And Name | Extension | Publish URL | zh-cn | de-de | Latest HO Xliff Generate Date
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $fileName
$wsOverview.Range("B3").Value = $pathAndName
$wsOverview.Range("C3").Value = $extension
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = $status
$wsOverview.Range("F3").Value = $status
$wsOverview.Range("G3").Value = $hoDate

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $baseUrl, "", "", $pathAndName) | Out-Null

# ---------------------------------------------------------------
# Sheet "zh-cn": Source File Name | File Extension | Status | Source Path | Priority |
#   Content Duplicate | Latest Handoff File | Latest Handoff Datetime | Latest Target File |
#   Latest Handback File | Latest Handback DateTime | Reference Tokens | To be localized |
#   Dependency From | Has metadata | Error Detail
# ---------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

$wsZh.Range("A3").Value = $fileName
$wsZh.Range("B3").Value = $extension
$wsZh.Range("C3").Value = $status
$wsZh.Range("D3").Value = $sourcePath
$wsZh.Range("E3").Value = $priority
$wsZh.Range("F3").Value = $contentDup
$wsZh.Range("G3").Value = $zhXlf
$wsZh.Range("H3").Value = $hoDateZh
$wsZh.Range("I3").Value = ""
$wsZh.Range("J3").Value = ""
$wsZh.Range("K3").Value = $handbackDt
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = $toLocalize
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = $contentDup
$wsZh.Range("P3").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $baseUrl, "", "", $fileName) | Out-Null

# ---------------------------------------------------------------
# Sheet "de-de": same column layout as zh-cn
# ---------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

$wsDe.Range("A3").Value = $fileName
$wsDe.Range("B3").Value = $extension
$wsDe.Range("C3").Value = $status
$wsDe.Range("D3").Value = $sourcePath
$wsDe.Range("E3").Value = $priority
$wsDe.Range("F3").Value = $contentDup
$wsDe.Range("G3").Value = $deXlf
$wsDe.Range("H3").Value = $hoDate
$wsDe.Range("I3").Value = ""
$wsDe.Range("J3").Value = ""
$wsDe.Range("K3").Value = $handbackDt
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = $toLocalize
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = $contentDup
$wsDe.Range("P3").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $baseUrl, "", "", $fileName) | Out-Null

# ---------------------------------------------------------------
# Column width auto-fit for the Status/zh-cn/de-de columns that now
# contain the longer "Ready for handoff" text.
# ---------------------------------------------------------------
$wsOverview.Columns.Item(5).AutoFit() | Out-Null
$wsOverview.Columns.Item(6).AutoFit() | Out-Null
$wsZh.Columns.Item(3).AutoFit() | Out-Null
$wsDe.Columns.Item(3).AutoFit() | Out-Null

Write-Host "Report rows for handoff added."
